$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -13.01
$ws.Range("C4").Value = -13.092
$ws.Range("D6").Value = -7.882000000000001
$ws.Range("C7").Value = -13.138
$ws.Range("D7").Value = -7.334999999999999
$ws.Range("C8").Value = -12.672
$ws.Range("D8").Value = -7.858
$ws.Range("B11").Value = 5.642000000000001
$ws.Range("E11").Value = 12.842
$ws.Range("B12").Value = 5.568
$ws.Range("C12").Value = -13.213
$ws.Range("C14").Value = -11.851
$ws.Range("E14").Value = 12.915
$ws.Range("B15").Value = 6.844999999999999
$ws.Range("D19").Value = -7.822
$ws.Range("E19").Value = 12.846
$ws.Range("D21").Value = -7.858
$ws.Range("E21").Value = 13.378
$ws.Range("C22").Value = -13.318
$ws.Range("D24").Value = -7.934
$ws.Range("D25").Value = -7.861
